$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update company/contact record in column B: "JORDI" -> "D' SANTI" ---
# (Assignment order matters for shared-string table allocation order.)
$ws.Range("B1").Value = "D' SANTI"
$ws.Range("B4").Value = "2445-8293"
$ws.Range("B5").Value = "200m Sur Urgencias Hospital San Ramon"

# --- Rebuild hyperlinks: B8 keeps rId1 (new email target), B7 gains a new hyperlink (rId2) ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B8"), "mailto:dsanti.srm@gmail.com")
$ws.Range("B8").Value = "dsanti.srm@gmail.com"
$ws.Range("B8").Style = "Hipervínculo"

$ws.Range("B9").Value = "Srmdsanti2"

$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:mariela.arcef@gmail.com")
$ws.Range("B7").Style = "Hipervínculo"

# --- Restore recorded window scroll position for the workbook view ---
$excel.ActiveWindow.Top = 912
